$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 115 (pushes the existing rows 115-133 down to 116-134)
$ws.Rows.Item(115).Insert()

# Populate the new derived-variable row: Rx16 / plasma / Treatments / description
$ws.Cells.Item(115, 1).Value = "Rx16"
$ws.Cells.Item(115, 2).Value = "plasma"
$ws.Cells.Item(115, 3).Value = "Treatments"
$ws.Cells.Item(115, 4).Value = "Convalescent plasma as COVID-19 treatment ever"

# Grow Table1 so it covers the newly inserted row (A1:E133 -> A1:E134)
$lo = $ws.ListObjects.Item(1)
$lastRow = $lo.Range.Rows.Count + 1
$lo.Resize($ws.Range("A1:E" + $lastRow))

# Match the saved selection / active cell from the source edit
$ws.Range("D115").Select()
